$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Week"
$ws.Range("B1").Value = "Predicted_Quantity"

# Remove the old Predicted_Quantity column (C) entirely - delete the column
$ws.Range("C1:C9").Delete()

# Update data rows: Week labels (A2:A9) and Predicted_Quantity values (B2:B9)
$weeks = @("2025-W43", "2025-W44", "2025-W45", "2025-W46", "2025-W47", "2025-W48", "2025-W49", "2025-W50")
$values = @(19, 15, 19, 16, 11, 17, 16, 13)

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $weeks[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}
